$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("lifts")

$data = @(
    @(6,  "NONE", "047G", "047G-P", "Olispa munkkeja"),
    @(7,  "NONE", "047G", "047G-P", "LäähPuuh"),
    @(8,  "NONE", "047G", "047G-P", "Note"),
    @(9,  "NONE", "S",    "O",      "LäähPuuh"),
    @(10, "NONE", "S",    "O",      "Note"),
    @(11, "NONE", "S",    "O",      "Note"),
    @(12, "NONE", "047G", "047G-P", "Ryys"),
    @(13, "NONE", "047G", "047G-P", "Note"),
    @(14, "NONE", "047G", "O",      "Konde")
)

$startRow = 9
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    $ws.Cells.Item($row, 1).Value = $values[0]
    $ws.Cells.Item($row, 2).Value = $values[1]
    $ws.Cells.Item($row, 3).Value = $values[2]
    $ws.Cells.Item($row, 4).Value = $values[3]
    $ws.Cells.Item($row, 5).Value = $values[4]
}
